$wb = $excel.ActiveWorkbook

function Set-SheetValues {
    param($ws, $values, $hyperlinkRefs)

    # Build a lookup of existing hyperlinks on this sheet, keyed by their
    # cell address (e.g. "$A$2"), so we can update the visible display
    # text without disturbing the underlying target (r:id / URL).
    $hlMap = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        $hlMap[$addr] = $hl
    }

    foreach ($ref in $values.Keys) {
        $newValue = $values[$ref]
        $ws.Range($ref).Value = $newValue

        if ($hyperlinkRefs -contains $ref) {
            $addrKey = "`$" + $ref.Substring(0,1) + "`$" + $ref.Substring(1)
            if ($hlMap.ContainsKey($addrKey)) {
                $hlMap[$addrKey].TextToDisplay = $newValue
            }
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": the two source files (5bbc677a... and c8d950fd...)
# swap rows, and the 5bbc677a row is now "Ready for handoff" with a new
# handoff date, while c8d950fd remains "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$overviewValues = @{
    "A2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md"
    "A3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md"
    "B3" = "Ready for handoff"
    "C3" = "Ready for handoff"
    "D3" = "2016-38-12 04:38:05"
}
$overviewHyperlinks = @("A2", "A3")
Set-SheetValues $wsOverview $overviewValues $overviewHyperlinks

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)
$zhValues = @{
    "A2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md"
    "D2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.zh-cn.xlf"
    "F2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md"
    "G2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.zh-cn.xlf"
    "A3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md"
    "C3" = "Ready for handoff"
    "D3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.zh-cn.xlf"
    "E3" = "2016-03-12 04:38:02"
    "F3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md"
    "G3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.zh-cn.xlf"
}
$zhHyperlinks = @("A2", "B2", "D2", "F2", "G2", "A3", "B3", "D3", "F3", "G3")
Set-SheetValues $wsZh $zhValues $zhHyperlinks

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)
$deValues = @{
    "A2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md"
    "D2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.de-de.xlf"
    "F2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md"
    "G2" = "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.de-de.xlf"
    "A3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md"
    "C3" = "Ready for handoff"
    "D3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.de-de.xlf"
    "E3" = "2016-03-12 04:38:05"
    "F3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md"
    "G3" = "5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.de-de.xlf"
}
$deHyperlinks = @("A2", "B2", "D2", "F2", "G2", "A3", "B3", "D3", "F3", "G3")
Set-SheetValues $wsDe $deValues $deHyperlinks
